$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update calibration / demand values on row 2 (Malawi specs)
$ws.Range("G2").Value = 2858.003227632393      # UrbanCutOff
$ws.Range("H2").Value = 0.1598539884266991     # UrbanRatioModelled
$ws.Range("V2").Value = 0.09908597401481424    # ElecModelled
$ws.Range("W2").Value = 10                     # MinNightLights
$ws.Range("Y2").Value = 1                      # MaxGridDist
$ws.Range("AA2").Value = 6140.790349047872     # PopCutOffRoundOne
$ws.Range("AB2").Value = 6000                  # PopCutOffRoundTwo
